$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "285.85"
Set-TextValue "E2" "2.44%"
Set-TextValue "E3" "4.73%"
Set-TextValue "D4" "5.044"
Set-TextValue "E4" "4.24%"
Set-TextValue "D5" "0.06695"
Set-TextValue "E5" "4.99%"
Set-TextValue "D6" "7.346"
Set-TextValue "E6" "4.46%"
Set-TextValue "D7" "3.386"
Set-TextValue "E7" "1.79%"
Set-TextValue "D8" "1.369"
Set-TextValue "E8" "5.30%"
Set-TextValue "D9" "0.9404"
Set-TextValue "D10" "0.1568"
Set-TextValue "E10" "2.95%"
Set-TextValue "D11" "0.06727"
Set-TextValue "E11" "13.80%"
Set-TextValue "D12" "0.07591"
Set-TextValue "E12" "1.00%"
Set-TextValue "D13" "0.02946"
Set-TextValue "E13" "0.95%"
Set-TextValue "D14" "0.09001"
Set-TextValue "E14" "0.05%"
Set-TextValue "D15" "0.001599"
Set-TextValue "E15" "2.00%"
Set-TextValue "D16" "0.04491"
Set-TextValue "E16" "1.64%"
Set-TextValue "D17" "0.0006455"
Set-TextValue "E17" "1.04%"
Set-TextValue "D18" "0.006659"
Set-TextValue "E18" "9.30%"
Set-TextValue "D19" "3.448"
Set-TextValue "E19" "-0.91%"
Set-TextValue "E20" "1.03%"
Set-TextValue "D21" "0.3209"
Set-TextValue "E21" "1.96%"
Set-TextValue "D22" "0.1310"
Set-TextValue "E22" "-3.00%"
Set-TextValue "D23" "4.076"
Set-TextValue "E23" "4.25%"
Set-TextValue "D24" "0.1549"
Set-TextValue "E24" "3.09%"
Set-TextValue "D25" "0.001179"
Set-TextValue "E25" "0.31%"
Set-TextValue "D26" "0.004491"
Set-TextValue "E26" "4.93%"
Set-TextValue "D27" "0.0001247"
Set-TextValue "E27" "5.66%"
Set-TextValue "D28" "0.0001614"
Set-TextValue "E28" "-2.36%"
Set-TextValue "D40" "0.04195"
Set-TextValue "E40" "3.00%"
Set-TextValue "D41" "0.006735"
Set-TextValue "E41" "1.63%"
Set-TextValue "D42" "0.1253"
Set-TextValue "E42" "-11.03%"
Set-TextValue "D43" "0.002014"
Set-TextValue "E43" "-3.60%"
Set-TextValue "D44" "0.01229"
Set-TextValue "E44" "11.72%"
Set-TextValue "D45" "0.00005658"
Set-TextValue "E45" "2.33%"
Set-TextValue "E46" "25.93%"
Set-TextValue "D47" "0.01304"
Set-TextValue "E47" "-29.48%"
